# Final outcome measures added: fill in the "Post Treatment" (column D)
# results for each question, matching the pattern already used in the
# Pre Baseline / Pre Experimental columns (B and C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "A little stressful"
$ws.Range("D3").Value = "Moderately stressful"
$ws.Range("D4").Value = "Moderately stressful"
$ws.Range("D5").Value = "Not stressful"
$ws.Range("D6").Value = "A little stressful"
$ws.Range("D7").Value = "Moderately stressful"

# Widen column D to match column B now that it holds text answers too.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Leave the active cell where data entry finished, one row below the table.
$ws.Range("D8").Select()
